# Rework the five "challenge" bullet points under "Expansión corporativa
# planificada": move the colon so it follows the full bold heading phrase
# instead of splitting mid-phrase, and tweak a few words in the bodies.
#
# Each bullet is two runs: a bold "heading" run and a normal "body" run.
# To keep the bold/non-bold split intact we replace each run's text with
# its own Find/Replace call (never spanning the run boundary in one call).

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "NOT FOUND: $old"
    }
}

# 1. Reconocimiento / Concienciación y reconocimiento de marca limitados
Replace-Text "Reconocimiento y reconocimiento" "Concienciación y reconocimiento de marca limitados"
Replace-Text " de marca limitados: lograr visibilidad en estos nuevos mercados es un obstáculo principal, lo que requiere un sólido esfuerzo de marketing para crear la presencia de marca de Adatum desde cero." ": un obstáculo principal es lograr visibilidad en estos nuevos mercados, lo que requiere un sólido esfuerzo de marketing para crear la presencia de marca de Adatum desde cero."

# 2. Intensa competencia / Competencia intensa
Replace-Text "Intensa competencia" "Competencia intensa"
Replace-Text ": el sector de servicios en la nube en Canadá es ferozmente competitivo, con numerosos jugadores." ": el sector de servicios en la nube en Canadá es ferozmente competitivo, con numerosos participantes."

# 3. Diversas preferencias y expectativas / Preferencias y expectativas diversas de los clientes
Replace-Text "Diversas preferencias y expectativas" "Preferencias y expectativas diversas de los clientes"
Replace-Text " de los clientes: adaptar productos y marketing para alinearse con las diversas demandas de estos mercados es fundamental para replicar con empresas y consumidores locales." ": adaptar los productos y el marketing para alinearse con las diversas demandas de estos mercados es fundamental para resonar con las empresas y los consumidores locales."

# 4. Desafíos / Desafíos normativos y de cumplimiento
Replace-Text "Desafíos" "Desafíos normativos y de cumplimiento"
Replace-Text " normativos y de cumplimiento: Adatum se enfrenta a la compleja tarea de navegar por la privacidad, la seguridad y las regulaciones operativas de la región, lo que necesita esfuerzos de cumplimiento diligentes." ": Adatum se enfrenta a la compleja tarea de navegar por la privacidad, la seguridad y las regulaciones operativas de la región, requiriendo esfuerzos de cumplimiento diligentes."

# 5. Complejidad / Complejidad operativa y logística
Replace-Text "Complejidad" "Complejidad operativa y logística"
Replace-Text " operativa y logística: el establecimiento de operaciones eficientes entre regiones presenta desafíos logísticos, especialmente en el mantenimiento de altos niveles de servicio y la administración de centros de datos en ubicaciones geográficas." ": el establecimiento de operaciones eficientes entre regiones presenta desafíos logísticos, especialmente en el mantenimiento de altos niveles de servicio y la administración de centros de datos en todas las ubicaciones geográficas."

Write-Output "Done"
